$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.728.64"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.82%  "

$ws.Range("D3").Value = "'3.166.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.77%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'590.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.52%  "

$ws.Range("D6").Value = "'133.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.35%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'3.164.99"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.78%  "

$ws.Range("E9").Value = "  -0.66%  "

$ws.Range("E10").Value = "  -6.93%  "

$ws.Range("D11").Value = "'5.22"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.20%  "

$ws.Range("E12").Value = "  -3.98%  "

$ws.Range("D13").Value = "'0.0000235"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.41%  "

$ws.Range("D14").Value = "'34.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.16%  "

$ws.Range("D15").Value = "'3.689.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.76%  "

$ws.Range("E16").Value = "  -1.14%  "

$ws.Range("D17").Value = "'3.165.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.92%  "

$ws.Range("D18").Value = "'62.747.20"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.94%  "

$ws.Range("E19").Value = "  -4.94%  "

$ws.Range("D20").Value = "'459.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.75%  "

$ws.Range("D21").Value = "'13.95"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.10%  "

$ws.Range("E22").Value = "  -6.45%  "

$ws.Range("D23").Value = "'7.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.70%  "

$ws.Range("D24").Value = "'13.36"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.49%  "

$ws.Range("D25").Value = "'82.81"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.56%  "

$ws.Range("E26").Value = "  +0.00%  "

$ws.Range("E27").Value = "  -0.06%  "

$ws.Range("D28").Value = "'2.66"
$ws.Range("D28").Style = "Normal"

$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.40%  "

$ws.Range("B30").Value = "NEARProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D30").Value = "'6.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.27%  "

$ws.Range("E31").Value = "  -6.61%  "

$ws.Range("D32").Value = "'27.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.36%  "

$ws.Range("E33").Value = "  -4.26%  "

$ws.Range("D34").Value = "'2.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.11%  "

$ws.Range("E35").Value = "  -6.76%  "

$ws.Range("D36").Value = "'5.79"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.93%  "

$ws.Range("D37").Value = "'51.21"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.32%  "

$ws.Range("E38").Value = "  -6.91%  "

$ws.Range("D39").Value = "'0.0386"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.54%  "

$ws.Range("D40").Value = "'403.21"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.42%  "

$ws.Range("D41").Value = "'8.04"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.00%  "

$ws.Range("E42").Value = "  -5.78%  "

$ws.Range("D43").Value = "'2.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.89%  "

$ws.Range("D44").Value = "'2.788.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.97%  "

$ws.Range("E45").Value = "  -6.77%  "

$ws.Range("E47").Value = "  -7.03%  "

$ws.Range("D48").Value = "'124.44"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.24%  "

$ws.Range("D49").Value = "'25.13"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.10%  "

$ws.Range("D50").Value = "'34.26"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.20%  "

$ws.Range("E51").Value = "  -2.39%  "
